# Commit: feat: add 2022-Q4 data
#
# This workbook tracks a stock's holders (funds) broken down by quarter.
# A new quarter "2022-Q4" is being added:
#   1. A brand-new worksheet "2022-Q4" (holding the detailed fund list for
#      that quarter) is inserted right after "总计" and before "2022-Q2".
#   2. The "总计" (summary) worksheet gets a new row for 2022-Q4 at the top
#      of its data table (the other quarters all stay the same, just shifted
#      down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Final state of the summary table (date label, holding count, market value).
$summaryRows = @(
    @("2022-Q4", 5,  0.13),
    @("2022-Q2", 3,  0.08),
    @("2022-Q1", 2,  1.66),
    @("2021-Q4", 2,  0.23),
    @("2021-Q3", 1,  0.38),
    @("2021-Q2", 2,  0.78),
    @("2021-Q1", 12, 0.14),
    @("2020-Q4", 7,  2.81)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]

    $aCell = $summary.Cells.Item($r, 1)
    $aCell.Value = $i                              # column A - running index (0-based)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108              # xlCenter
    $aCell.VerticalAlignment = -4160                # xlTop
    $aCell.Borders.LineStyle = 1

    $summary.Cells.Item($r, 2).Value = $row[0]     # column B - quarter label
    $summary.Cells.Item($r, 3).Value = $row[1]     # column C - holding count
    $summary.Cells.Item($r, 4).Value = $row[2]     # column D - market value (亿元)
}

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q4.Cells.Item(1, $c + 2)   # headers start at column B
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

# Fund detail rows. Columns B-G are stored as text (to preserve fund codes'
# leading zeros and match the source formatting), A and H are numeric.
$fundRows = @(
    @("010568", "海富通惠睿精选混合A",     "11.95", "29.51", "0.58", "0.0693", 8),
    @("010569", "海富通惠睿精选混合C",     "5.11",  "29.51", "0.58", "0.0296", 8),
    @("004703", "南方兴盛先锋灵活配置混合", "0.74",  "63.23", "3.72", "0.0275", 5),
    @("006818", "安信盈利驱动股票A",       "0.08",  "88.09", "4.40", "0.0035", 8),
    @("006819", "安信盈利驱动股票C",       "0.06",  "88.09", "4.40", "0.0026", 8)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $aCell = $q4.Cells.Item($r, 1)
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    # Force text storage for columns B-G so numeric-looking strings (fund
    # codes, percentages, ...) keep their original textual representation.
    $textRange = $q4.Range($q4.Cells.Item($r, 2), $q4.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"

    $q4.Cells.Item($r, 2).Value = $row[0]   # 基金代码
    $q4.Cells.Item($r, 3).Value = $row[1]   # 基金名称
    $q4.Cells.Item($r, 4).Value = $row[2]   # 基金规模
    $q4.Cells.Item($r, 5).Value = $row[3]   # 股票总仓位
    $q4.Cells.Item($r, 6).Value = $row[4]   # 仓位占比
    $q4.Cells.Item($r, 7).Value = $row[5]   # 持有市值(亿元)
    $q4.Cells.Item($r, 8).Value = $row[6]   # 仓位排名 (numeric)
}

# Keep the originally-active "2020-Q4" tab selected (adding a sheet would
# otherwise leave the brand-new "2022-Q4" sheet focused).
$wb.Worksheets.Item("2020-Q4").Select()

